$wb = $excel.ActiveWorkbook

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.0835791708173058
$wsLinear.Range("B3").Value = 3.677798340639166
$wsLinear.Range("B4").Value = 924.5770670424745
$wsLinear.Range("B5").Value = "[1.0, 0.2219175969280435, 0.07820482702373979, 0.12934056071917785, 0.08718744929392527, 0.04811436564015423, 0.1956915170073739, 0.3367260556058545, 0.17611139220178573, 0.035925830658739964, 0.07052862951986745, 0.0810116737968763, 0.02254840769355749, 0.16594972618090406, 0.32132208397379625, 0.15022154421323144, 0.00275541430796133, 0.05809578029869869, 0.05011030506278425, -0.00842494278007941]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = 0.07810602730467787
$wsNonLinear.Range("B5").Value = 4.434734591741209
$wsNonLinear.Range("B6").Value = 954.3715624016675
$wsNonLinear.Range("B7").Value = 0.7716092784002508
$wsNonLinear.Range("B8").Value = 1.993037058717917
$wsNonLinear.Range("B9").Value = 899.5814637093338
$wsNonLinear.Range("B10").Value = "[0.9999999999999999, 0.2213466242475205, 0.07888934827472299, 0.13056672309671116, 0.08869891259808177, 0.04893962242323496, 0.1955436243557061, 0.3355776206995872, 0.17590937937808987, 0.03646401684858577, 0.07177830521233783, 0.08222912365208797, 0.02324596276383912, 0.16587603700564602, 0.3206657637269602, 0.15051102080368886, 0.0038222597478946314, 0.05888788940559675, 0.05076867904840815, -0.007960980997903083]"
